$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $text
    $rng.Style = "Normal"
}

# Row 2 (Bitcoin)
Set-TextValue "D2" "67.122.48"
$ws.Range("E2").Value = "  -0.79%  "

# Row 3 (Ethereum)
Set-TextValue "D3" "2.610.99"

# Row 4 (TetherUSD)
$ws.Range("E4").Value = "  +0.15%  "

# Row 5 (BNB)
Set-TextValue "D5" "593.77"
$ws.Range("E5").Value = "  -0.34%  "

# Row 6 (Solana)
Set-TextValue "D6" "152.30"
$ws.Range("E6").Value = "  -2.37%  "

# Row 8 (XRP)
Set-TextValue "D8" "0.557"
$ws.Range("E8").Value = "  +2.55%  "

# Row 9 (LidoStakedEther)
Set-TextValue "D9" "2.608.99"
$ws.Range("E9").Value = "  -0.26%  "

# Row 11 (TRON)
$ws.Range("E11").Value = "  +0.33%  "

# Row 12 (Toncoin)
Set-TextValue "D12" "5.15"

# Row 13 (Cardano)
Set-TextValue "D13" "0.345"
$ws.Range("E13").Value = "  -3.30%  "

# Row 14 (Avalanche)
Set-TextValue "D14" "27.49"
$ws.Range("E14").Value = "  -0.74%  "

# Row 15 (WrappedliquidstakedEther2.0)
Set-TextValue "D15" "3.085.63"
$ws.Range("E15").Value = "  -0.01%  "

# Row 16 (ShibaInu)
$ws.Range("E16").Value = "  -3.95%  "

# Row 17 (WrappedBTC)
Set-TextValue "D17" "66.966.78"
$ws.Range("E17").Value = "  -0.50%  "

# Row 18 (WrappedEther)
Set-TextValue "D18" "2.608.32"
$ws.Range("E18").Value = "  -0.60%  "

# Row 19 (BitcoinCash)
Set-TextValue "D19" "365.50"
$ws.Range("E19").Value = "  +0.97%  "

# Row 20 (Chainlink)
Set-TextValue "D20" "10.99"
$ws.Range("E20").Value = "  -3.30%  "

# Row 21 (Uniswap)
$ws.Range("E21").Value = "  -5.27%  "

# Row 22 (Polkadot)
$ws.Range("E22").Value = "  -0.62%  "

# Row 23 (SuiNetwork)
$ws.Range("E23").Value = "  -1.05%  "

# Row 24 (Dai)
Set-TextValue "D24" "0.999"
$ws.Range("E24").Value = "  -0.14%  "

# Row 25 (Aptos)
Set-TextValue "D25" "10.01"
$ws.Range("E25").Value = "  -0.81%  "

# Row 26 (Litecoin)
Set-TextValue "D26" "66.63"
$ws.Range("E26").Value = "  -6.85%  "

# Row 27 (WrappedeETH)
Set-TextValue "D27" "2.747.07"
$ws.Range("E27").Value = "  -0.42%  "

# Row 28 and 29 swap places: Bittensor now at 28, Binance-PegBSC-USD now at 29
$ws.Range("B28").Value = "Bittensor"
$ws.Range("C28").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D28" "581.19"
$ws.Range("E28").Value = "  -2.39%  "

$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue "D29" "0.999"
$ws.Range("E29").Value = "  -0.67%  "

# Row 30 (PEPE)
Set-TextValue "D30" "0.0000101"
$ws.Range("E30").Value = "  -3.58%  "

# Row 31 (Fetch.AI)
$ws.Range("E31").Value = "  -4.75%  "

# Row 32 (InternetComputer(DFINITY))
Set-TextValue "D32" "7.74"
$ws.Range("E32").Value = "  -3.01%  "

# Row 33 (PancakeSwap)
$ws.Range("E33").Value = "  -2.24%  "

# Row 34 (FirstDigitalUSD)
$ws.Range("E34").Value = "  +0.02%  "

# Row 35 (Kaspa)
Set-TextValue "D35" "0.122"
$ws.Range("E35").Value = "  -8.79%  "

# Row 36 (ImmutableX)
$ws.Range("E36").Value = "  -3.44%  "

# Row 37 (NEARProtocol)
Set-TextValue "D37" "4.86"
$ws.Range("E37").Value = "  -2.37%  "

# Row 38 (Monero)
Set-TextValue "D38" "156.38"
$ws.Range("E38").Value = "  +1.25%  "

# Row 39 (EthereumClassic)
Set-TextValue "D39" "18.98"
$ws.Range("E39").Value = "  -3.18%  "

# Row 40 (PolygonEcosystemToken)
Set-TextValue "D40" "0.366"
$ws.Range("E40").Value = "  -1.63%  "

# Row 41 (RenderToken)
$ws.Range("E41").Value = "  -3.64%  "

# Row 42 (Stacks)
Set-TextValue "D42" "1.80"
$ws.Range("E42").Value = "  -2.83%  "

# Row 43 (dogwifhat)
$ws.Range("E43").Value = "  -2.78%  "

# Row 44 (OKB)
Set-TextValue "D44" "41.08"
$ws.Range("E44").Value = "  -0.67%  "

# Row 45 (USDe)
$ws.Range("E45").Value = "  -0.01%  "

# Row 46 (WhiteBITCoin)
Set-TextValue "D46" "16.39"
$ws.Range("E46").Value = "  -0.61%  "

# Row 47 (Aave)
Set-TextValue "D47" "155.77"
$ws.Range("E47").Value = "  -1.48%  "

# Row 48 (BabyDogeCoin)
Set-TextValue "D48" "0.0₆0286"
$ws.Range("E48").Value = "  -3.16%  "

# Row 49 (Filecoin)
$ws.Range("E49").Value = "  -0.75%  "

# Row 50 (InjectiveProtocol)
Set-TextValue "D50" "21.74"
$ws.Range("E50").Value = "  +3.67%  "

# Row 51 (Mantle)
$ws.Range("E51").Value = "  -2.23%  "
